$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing "Florian Kleinig" contact row (row 5) to "Jack  Kleinig" ---
# Set E5 before D5 so the shared-string table gets the same ordering as the target file
# (23 = "Jack  Kleinig", 24 = "jack.kleinig@droniq.de").
$ws.Range("E5").Value = "Jack  Kleinig"
$ws.Range("D5").Value = "jack.kleinig@droniq.de"

# --- Add two new rows for Drone Volt contacts ---
# Row 6: Olivier Gualdoni
$ws.Range("A6").Value = "yes"
$ws.Range("B6").Value = "Drone Volt "
$ws.Range("C6").Value = "https://www.dronevolt.com/`n"
$ws.Range("E6").Value = "Olivier Gualdoni"
$ws.Range("D6").Value = "olivier@dronevolt.com"

# Row 7: Stefano Valentini
$ws.Range("A7").Value = "yes"
$ws.Range("B7").Value = "Drone Volt "
$ws.Range("C7").Value = "https://www.dronevolt.com/`n"
$ws.Range("E7").Value = "Stefano Valentini"
$ws.Range("D7").Value = "stefano.valentini@dronevolt.com"

# The multi-line URL text triggers an automatic row-height change; re-run AutoFit
# so the rows fall back to the sheet's default height (no explicit ht/customHeight),
# matching the target layout.
$ws.Rows(6).AutoFit()
$ws.Rows(7).AutoFit()

# --- Hyperlinks ---
# The "Florian Kleinig" mailto hyperlink on D5 must now point at Jack's address, but
# this COM layer only supports adding hyperlinks (re-assigning .Address on an existing
# hyperlink duplicates it instead of replacing it in place). So drop every hyperlink on
# the sheet and recreate them all, in their original order, pointing at the right
# addresses, then finish with the four brand-new ones for rows 6-7.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.google.com/")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:testing1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://droniq.de/en/?")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:alexander.wulf@droniq.de")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://droniq.de/en/?")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.google.com/")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:testing22@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:jack.kleinig@droniq.de")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.dronevolt.com/")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:olivier@dronevolt.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:stefano.valentini@dronevolt.com")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://www.dronevolt.com/")

# Re-apply the built-in Hyperlink style only to the cells that used it originally
# (C3, C4, D4, C5, D5) plus the new hyperlink cells, so style indices match the
# target - C2/D2/D3 never had the Hyperlink style applied in the source file, so put
# those back to Normal (Hyperlinks.Add auto-applies the Hyperlink style to them too).
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("C6").Style = "Hyperlink"
$ws.Range("D6").Style = "Hyperlink"
$ws.Range("C7").Style = "Hyperlink"
$ws.Range("D7").Style = "Hyperlink"

# --- Column D width (19 -> 22.125 character units) ---
$ws.Columns("D").ColumnWidth = 21.505

# --- Selection moved to E14 ---
$ws.Range("E14").Select()
